# Edit script for Poster_Unlinkability.pptx
# Commit message: "Adding new security risk"
#
# Shape index map (1-based, Slide 1), discovered via inspection:
#   4  -> id=9  "Textfeld 8"   (Summary / Introduction column)
#   5  -> id=11 "Textfeld 10"  (.../Goal column)
#   6  -> id=18 "Textfeld 17"  (Results column)
#   7  -> id=19 "Picture 18"   (OIDC4VP / trust-triangle picture)
#   8  -> id=21 "Textfeld 4"   ("The Trust triangle" caption)
#   9  -> id=4  "Picture 3"    (example-VC picture)
#   10 -> id=5  "Textfeld 4"   ("Example of a VC" caption)

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Shape 4 (id=9, "Textfeld 8") - Summary / Introduction text column
# ---------------------------------------------------------------------------
$sh9 = $s.Shapes.Item(4)
$sh9.Height = 16866156 / $EMU_PER_PT

$tr9 = $sh9.TextFrame.TextRange
$full9 = $tr9.Text

# "Summary" (first heading) -> "Introduction"
$idx = $full9.IndexOf("Summary")
$rng = $tr9.Characters($idx + 1, "Summary".Length)
$rng.Text = "Introduction"

# refresh text after edit, then update the second heading
$full9 = $tr9.Text
$idx = $full9.IndexOf("Introduction", $idx + "Introduction".Length)
$rng = $tr9.Characters($idx + 1, "Introduction".Length)
$rng.Text = "SSI and VCs"

# Extend the last paragraph: replace the trailing
#   " (BBS) is used in this thesis. In physical credentials there are"
# with ", created by Dan " plus several runs naming the BBS authors,
# ending in " (BBS) is" (the rest moves to the next text box).
$full9 = $tr9.Text
$oldTail = " (BBS) is used in this thesis. In physical credentials there are"
$idx = $full9.IndexOf($oldTail)
$rng = $tr9.Characters($idx + 1, $oldTail.Length)
$rng.Text = ", created by Dan "

$tail9 = $tr9.Characters($tr9.Length, 1)
$tail9.InsertAfter("Boneh") | Out-Null
$tail9 = $tr9.Characters($tr9.Length, 1)
$tail9.InsertAfter(", Xavier ") | Out-Null
$tail9 = $tr9.Characters($tr9.Length, 1)
$tail9.InsertAfter("Boyen") | Out-Null
$tail9 = $tr9.Characters($tr9.Length, 1)
$tail9.InsertAfter(", and ") | Out-Null
$tail9 = $tr9.Characters($tr9.Length, 1)
$tail9.InsertAfter("Hovav") | Out-Null
$tail9 = $tr9.Characters($tr9.Length, 1)
$tail9.InsertAfter(" ") | Out-Null
$tail9 = $tr9.Characters($tr9.Length, 1)
$tail9.InsertAfter("Shacham") | Out-Null
$tail9 = $tr9.Characters($tr9.Length, 1)
$tail9.InsertAfter(" (BBS) is") | Out-Null

# ---------------------------------------------------------------------------
# Shape 5 (id=11, "Textfeld 10") - continuation text column + Goal
# ---------------------------------------------------------------------------
$sh11 = $s.Shapes.Item(5)
$sh11.Top = 6482980 / $EMU_PER_PT
$sh11.Height = 11341566 / $EMU_PER_PT

$tr11 = $sh11.TextFrame.TextRange

# Prepend the continuation of the sentence moved from the previous text box.
$first11 = $tr11.Characters(1, 1)
$first11.InsertBefore("used in this thesis. In physical credentials there are ") | Out-Null

# Update the "Goal" paragraph's closing sentence.
$full11 = $tr11.Text
$oldGoal = "The goal is the analysis of the different technologies working together in a real-world use case. "
$idx = $full11.IndexOf($oldGoal)
$rng = $tr11.Characters($idx + 1, $oldGoal.Length)
$rng.Text = "The goal of this thesis is to analyze if using these different technologies together in a real-world use case, breaks the "

$tail11 = $tr11.Characters($tr11.Length, 1)
$tail11.InsertAfter("unlinkabilty") | Out-Null
$tail11 = $tr11.Characters($tr11.Length, 1)
$tail11.InsertAfter(" provided by BBS.") | Out-Null

# ---------------------------------------------------------------------------
# Shape 6 (id=18, "Textfeld 17") - Results column
# ---------------------------------------------------------------------------
$sh18 = $s.Shapes.Item(6)
$sh18.Height = 7817525 / $EMU_PER_PT

$tr18 = $sh18.TextFrame.TextRange
$full18 = $tr18.Text
$oldResult = "Knowing how to generate digital credentials using VCs, protect them using BBS and deliver them trough secure channels using OIDC4VP, shows a future where SSI is possible. Using the mentioned technologies as a basis, future research may contribute to a more secure digital world for individuals."
$idx = $full18.IndexOf($oldResult)
$rng = $tr18.Characters($idx + 1, $oldResult.Length)
$rng.Text = "The results of this thesis show, that using these technologies together, a future where SSI is the standard, is possible. Using the mentioned technologies as a basis, future research may contribute to a more secure digital world for individuals."

# ---------------------------------------------------------------------------
# Shape 7 (id=19, "Picture 18") - picture moved/resized
# ---------------------------------------------------------------------------
$sh19 = $s.Shapes.Item(7)
$sh19.Left = 11095507 / $EMU_PER_PT
$sh19.Width = 8088960 / $EMU_PER_PT
$sh19.Height = 4622263 / $EMU_PER_PT

# ---------------------------------------------------------------------------
# Shape 8 (id=21, "Textfeld 4") - "The Trust triangle" caption moved
# ---------------------------------------------------------------------------
$sh21 = $s.Shapes.Item(8)
$sh21.Left = 11027534 / $EMU_PER_PT
$sh21.Top = 5926581 / $EMU_PER_PT

# ---------------------------------------------------------------------------
# Shape 9 (id=4, "Picture 3") - picture moved/resized
# ---------------------------------------------------------------------------
$sh4 = $s.Shapes.Item(9)
$sh4.Left = 21116350 / $EMU_PER_PT
$sh4.Top = 9044043 / $EMU_PER_PT
$sh4.Width = 7547902 / $EMU_PER_PT
$sh4.Height = 8175677 / $EMU_PER_PT

# ---------------------------------------------------------------------------
# Shape 10 (id=5, "Textfeld 4") - "Example of a VC" caption moved/resized
# ---------------------------------------------------------------------------
$sh5 = $s.Shapes.Item(10)
$sh5.Left = 21116350 / $EMU_PER_PT
$sh5.Top = 17306145 / $EMU_PER_PT
$sh5.Width = 6642647 / $EMU_PER_PT
$sh5.Height = 523220 / $EMU_PER_PT
